$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected IFRS financial figures (error fix per commit "error solve ifrs list").
# Figures were previously inflated (appears an extra zero/scale or wrong source range
# was used); this replaces each numeric cell in rows 2-9 with the corrected value.
$rowData = @{
    2 = @{ "D" = 64661; "E" = 3140; "F" = 3140; "G" = 1916; "H" = 1127; "I" = 1127; "K" = 44288; "L" = 27323; "M" = 16965; "N" = 16965; "O" = 0; "P" = 1183; "Q" = 7430; "R" = -3088; "S" = -4402; "T" = 2609; "U" = 4821; "V" = 14603; "W" = 4.86; "X" = 1.74; "Y" = 7.48; "Z" = 2.5; "AA" = 161.05; "AB" = 1339.17; "AC" = 5020; "AD" = 22.41; "AE" = 71689; "AF" = 1.57; "AG" = 250; "AH" = 0.22; "AI" = 5.25; "AJ" = 23667107 }
    3 = @{ "D" = 61381; "E" = 2237; "F" = 2237; "G" = 1221; "H" = 951; "I" = 951; "K" = 39143; "L" = 21493; "M" = 17651; "N" = 17651; "O" = 0; "P" = 1183; "Q" = 6784; "R" = -3062; "S" = -4085; "T" = 2615; "U" = 4168; "V" = 10643; "W" = 3.64; "X" = 1.55; "Y" = 5.49; "Z" = 2.28; "AA" = 121.77; "AB" = 1399.15; "AC" = 4018; "AD" = 24.52; "AE" = 74586; "AF" = 1.32; "AG" = 350; "AH" = 0.36; "AI" = 8.710000000000001; "AJ" = 23667107 }
    4 = @{ "D" = 57546; "E" = 1048; "F" = 1048; "G" = 107; "H" = 50; "I" = 50; "K" = 43237; "L" = 25452; "M" = 17785; "N" = 17785; "O" = 0; "P" = 1183; "Q" = 3318; "R" = -3559; "S" = 65; "T" = 3190; "U" = 128; "V" = 10813; "W" = 1.82; "X" = 0.09; "Y" = 0.28; "Z" = 0.12; "AA" = 143.11; "AB" = 1416.42; "AC" = 209; "AD" = 422.41; "AE" = 75154; "AF" = 1.18; "AG" = 250; "AH" = 0.28; "AI" = 119.45; "AJ" = 23667107 }
    5 = @{ "D" = 76414; "E" = 2965; "F" = 2965; "G" = 2387; "H" = 1748; "I" = 1748; "K" = 58775; "L" = 39267; "M" = 19508; "N" = 19508; "O" = 0; "P" = 1183; "Q" = 4460; "R" = -8344; "S" = 4225; "T" = 7757; "U" = -3297; "V" = 15026; "W" = 3.88; "X" = 2.29; "Y" = 9.369999999999999; "Z" = 3.43; "AA" = 201.29; "AB" = 1584.48; "AC" = 7385; "AD" = 19.5; "AE" = 82435; "AF" = 1.75; "AG" = 250; "AH" = 0.17; "AI" = 3.38; "AJ" = 23667107 }
    6 = @{ "D" = 79821; "E" = 2635; "F" = 2635; "G" = 1840; "H" = 1631; "I" = 1631; "K" = 57561; "L" = 36378; "M" = 21182; "N" = 21182; "P" = 1183; "Q" = 10623; "R" = -14145; "S" = 6036; "T" = 14135; "U" = -3512; "V" = 21267; "W" = 3.3; "X" = 2.04; "Y" = 8.02; "Z" = 2.8; "AA" = 171.74; "AB" = 1724.22; "AC" = 6891; "AD" = 12.54; "AE" = 89511; "AF" = 0.97; "AG" = 300; "AH" = 0.35; "AI" = 4.35; "AJ" = 23667107 }
    7 = @{ "D" = 80819; "E" = 3666; "G" = 2271; "H" = 1830; "I" = 1806; "K" = 62177; "L" = 39205; "M" = 22972; "N" = 22983; "P" = 1181; "Q" = 7728; "R" = -6157; "S" = -1154; "T" = 6158; "U" = 1125; "W" = 4.54; "X" = 2.26; "Y" = 8.18; "Z" = 3.06; "AA" = 170.67; "AC" = 7633; "AD" = 19.98; "AE" = 97120; "AF" = 1.57; "AG" = 289; "AH" = 0.19; "AI" = 3.79 }
    8 = @{ "D" = 92932; "E" = 5226; "G" = 4312; "H" = 3331; "I" = 3332; "K" = 66621; "L" = 40994; "M" = 25627; "N" = 25649; "P" = 1181; "Q" = 9963; "R" = -7473; "S" = -762; "T" = 7071; "U" = 3084; "W" = 5.62; "X" = 3.59; "Y" = 13.69; "Z" = 5.17; "AA" = 159.96; "AC" = 14077; "AD" = 10.8; "AE" = 108386; "AG" = 324; "AH" = 0.21; "AI" = 2.3 }
    9 = @{ "D" = 99867; "E" = 5738; "G" = 4901; "H" = 3784; "I" = 3784; "K" = 71199; "L" = 41817; "M" = 29382; "N" = 29382; "P" = 1181; "Q" = 11248; "R" = -6741; "S" = -751; "T" = 6412; "U" = 4197; "W" = 5.75; "X" = 3.79; "Y" = 13.75; "Z" = 5.49; "AA" = 142.32; "AC" = 15988; "AD" = 9.51; "AE" = 124163; "AF" = 1.22; "AG" = 339; "AH" = 0.22; "AI" = 2.12 }
}

foreach ($r in $rowData.Keys) {
    foreach ($col in $rowData[$r].Keys) {
        $ws.Range("$col$r").Value = $rowData[$r][$col]
    }
}

# The per-row split of net income attributable to non-controlling interests (column J)
# no longer applies to the corrected figures in rows 2-5, so those cells are cleared
# (column J was already blank for rows 6-9).
foreach ($r in 2..5) {
    $ws.Range("J$r").ClearContents()
}
